# edit.ps1
# Applies the changes described by the target diff to the Word document:
#   1) Merge the two bold runs "git commit " + "-a -m "rotulo" " into a single run.
#   2) Merge the two plain runs "- Ao adicionar..." + "onde e possivel..." into a single run.
#   3) Append an empty paragraph followed by a bold+italic paragraph containing
#      "Teste de modificacao no arquivo".
#   4) Tweak a couple of paragraph-format properties on the "Normal" style
#      (reading order / alignment) that the target XML also shows changed.

$d = $word.ActiveDocument

function Merge-TwoRuns($paragraphIndex, $firstAnchor, $secondAnchor) {
    # Merges two adjacent runs inside a paragraph into one run, keeping the
    # formatting of the first run. $firstAnchor is text uniquely identifying
    # (from the start of) the first run, $secondAnchor uniquely identifies
    # (from the start of) the second run. Both anchors must lie inside the
    # given paragraph.
    $pRange = $d.Paragraphs($paragraphIndex).Range

    $f1 = $pRange.Duplicate
    $f1.Find.ClearFormatting()
    $f1.Find.Execute($firstAnchor, $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

    $f2 = $pRange.Duplicate
    $f2.Find.ClearFormatting()
    $f2.Find.Execute($secondAnchor, $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

    $splitPos = $f1.End
    $run2 = $d.Range($splitPos, $f2.End)
    $run2Text = $run2.Text
    $run2.Delete()

    $run1 = $d.Range($f1.Start, $splitPos)
    $run1.InsertAfter($run2Text)
}

# 1) "git commit " + "-a -m "rotulo" " -> merged into a single bold run.
Merge-TwoRuns 12 "git commit " "rótulo” "

# 2) "- Ao adicionar ... original, " + "onde e possivel ... especifica;" -> merged.
Merge-TwoRuns 18 "– Ao adicionar" "especifica;"

# 3) Append two new paragraphs at the end of the document.
$endRange = $d.Content
$endRange.Collapse(0)
$endRange.InsertParagraphAfter()

$endRange2 = $d.Content
$endRange2.Collapse(0)
$endRange2.InsertParagraphAfter()

$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$newTextRange = $lastPara.Range
$newTextRange.InsertAfter("Teste de modificação no arquivo")
$newTextRange.Font.Bold = 1
$newTextRange.Font.BoldBi = 1
$newTextRange.Font.Italic = 1
$newTextRange.Font.ItalicBi = 1

# 4) Normal style paragraph-format tweaks seen in the target styles.xml.
$normalStyle = $d.Styles("Normal")
$normalStyle.ParagraphFormat.ReadingOrder = 0
$normalStyle.ParagraphFormat.Alignment = 0
